$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking snapshot refresh (prices + 1h volume deltas).
# Cells whose new text is itself a plain decimal number need a leading
# apostrophe so Excel stores literal text (matching the source feed's
# formatting, e.g. trailing zeros) instead of silently recasting the cell
# as a Number; ".Style = 'Normal'" immediately after clears the resulting
# quote-prefix flag so the cell format stays plain, like the rest of the sheet.

# --- Row 30/31 swap: PEPE <-> Aptos (ranking reorder w/ refreshed data) ---
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").Value = "'6.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.67%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0723"
$ws.Range("E31").Value = "  +4.63%  "

# --- Price / Volume(1h) refresh for the rest of the table ---
$ws.Range("D2").Value = "56.614.52"
$ws.Range("E2").Value = "  +4.47%  "
$ws.Range("D3").Value = "2.318.83"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("D4").Value = "'0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "'517.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.18%  "
$ws.Range("D6").Value = "'133.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.63%  "
$ws.Range("D7").Value = "'0.995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("D8").Value = "'0.536"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.19%  "
$ws.Range("D9").Value = "2.343.70"
$ws.Range("E9").Value = "  +3.33%  "
$ws.Range("E10").Value = "  +8.43%  "
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("E12").Value = "  +7.98%  "
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("D14").Value = "'24.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.76%  "
$ws.Range("D15").Value = "2.733.71"
$ws.Range("E15").Value = "  +2.68%  "
$ws.Range("D16").Value = "56.389.83"
$ws.Range("E16").Value = "  +4.11%  "
$ws.Range("E17").Value = "  +4.46%  "
$ws.Range("D18").Value = "2.332.95"
$ws.Range("E18").Value = "  +3.17%  "
$ws.Range("D19").Value = "'10.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.78%  "
$ws.Range("E20").Value = "  +2.77%  "
$ws.Range("D21").Value = "'322.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.63%  "
$ws.Range("D22").Value = "'6.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.62%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "'61.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("E26").Value = "  +6.04%  "
$ws.Range("E27").Value = "  +4.77%  "
$ws.Range("D28").Value = "'171.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").Value = "'1.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.66%  "
$ws.Range("E32").Value = "  +4.50%  "
$ws.Range("D33").Value = "'18.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.49%  "
$ws.Range("D35").Value = "'0.997"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'1.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.00%  "
$ws.Range("D37").Value = "'0.927"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("D38").Value = "'3.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.38%  "
$ws.Range("E39").Value = "  +8.78%  "
$ws.Range("D40").Value = "'37.44"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.23%  "
$ws.Range("E41").Value = "  +2.59%  "
$ws.Range("E42").Value = "  +7.79%  "
$ws.Range("D43").Value = "'138.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.36%  "
$ws.Range("E44").Value = "  +7.25%  "
$ws.Range("D45").Value = "'268.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +11.34%  "
$ws.Range("D46").Value = "'0.0510"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.39%  "
$ws.Range("D48").Value = "'0.556"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("E49").Value = "  +2.75%  "
$ws.Range("E50").Value = "  +5.68%  "
$ws.Range("D51").Value = "'16.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.02%  "
